# Adds five new question/x-path row-pairs (rows 53-62) to Sheet1, mirroring
# the existing "label (A:E merged) | x-path (F:I merged)" layout used by the
# rest of the sheet, then blank spacer-style rows beneath each, and merges
# the appropriate cell ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108

$entries = @(
  @{ Row = 53; Label = "show more item";        Xpath = '// button[@class="_2KpZ6l _2jekE6" ]';                                                                      Wrap = $true  },
  @{ Row = 55; Label = "select particular item"; Xpath = '// span[@class="row _1kkfO3 BqOr_g" and text()="THE MAPLES FASHION Baby Girls Above Knee..." ]  ';          Wrap = $true  },
  @{ Row = 57; Label = "filter on the way";      Xpath = '// div[@class="YeQvMM" and text()="On the way"  ]  ';                                                       Wrap = $true  },
  @{ Row = 59; Label = "mannage address";        Xpath = '//div[@class="NS64GK" and text()="Manage Addresses"]';                                                      Wrap = $true  },
  @{ Row = 61; Label = "add new address";        Xpath = '//img[@class="_1g8pEu"]';                                                                                   Wrap = $false }
)

foreach ($entry in $entries) {
  $r1 = $entry.Row
  $r2 = $entry.Row + 1

  $labelBlock = $ws.Range("A" + $r1 + ":E" + $r2)
  $xpathBlock = $ws.Range("F" + $r1 + ":I" + $r2)

  $ws.Range("A" + $r1).Value = $entry.Label
  $ws.Range("F" + $r1).Value = $entry.Xpath

  $labelBlock.HorizontalAlignment = $xlCenter
  $xpathBlock.HorizontalAlignment = $xlCenter
  if ($entry.Wrap) {
    $xpathBlock.WrapText = $true
  }

  $labelBlock.Merge() | Out-Null
  $xpathBlock.Merge() | Out-Null
}

$ws.Range("F61:I62").Select() | Out-Null
